$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 4
$ws.Range("H2").Value = 2.08
$ws.Range("I2").Value = 2.16
$ws.Range("J2").Value = 3.35
$ws.Range("N2").Value = 3.35
$ws.Range("O2").Value = 1.37
$ws.Range("P2").Value = 1.81
$ws.Range("Q2").Value = 2.08
$ws.Range("R2").Value = 1.31
$ws.Range("S2").Value = 3.8
$ws.Range("T2").Value = 1.85
$ws.Range("V2").Value = 1.86
$ws.Range("W2").Value = 1.3
$ws.Range("X2").Value = 13
$ws.Range("Y2").Value = 9
$ws.Range("AA2").Value = 970
$ws.Range("AC2").Value = 9.2
$ws.Range("AD2").Value = 13
$ws.Range("AE2").Value = 29
$ws.Range("AH2").Value = 22
$ws.Range("AL2").Value = 65
$ws.Range("AN2").Value = 65
$ws.Range("AO2").Value = 22

# Row 3
$ws.Range("L3").Value = 1.01
$ws.Range("N3").Value = 2.02
$ws.Range("Q3").Value = 1.79
$ws.Range("S3").Value = 2.78
$ws.Range("T3").Value = 1.65
$ws.Range("U3").Value = 1.69
$ws.Range("Y3").Value = 38
$ws.Range("AB3").Value = 12
$ws.Range("AC3").Value = 15.5
$ws.Range("AD3").Value = 46
$ws.Range("AF3").Value = 13
$ws.Range("AG3").Value = 15
$ws.Range("AH3").Value = 36
$ws.Range("AJ3").Value = 19
$ws.Range("AK3").Value = 23

# Row 4
$ws.Range("F4").Value = 1.98
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 4.2
$ws.Range("I4").Value = 4.4
$ws.Range("J4").Value = 3.8
$ws.Range("K4").Value = 3.85
$ws.Range("N4").Value = 3.9
$ws.Range("O4").Value = 1.32
$ws.Range("P4").Value = 1.99
$ws.Range("Q4").Value = 1.97
$ws.Range("S4").Value = 3.45
$ws.Range("U4").Value = 2.14
$ws.Range("X4").Value = 15.5
$ws.Range("Y4").Value = 15.5
$ws.Range("Z4").Value = 32
$ws.Range("AA4").Value = 90
$ws.Range("AJ4").Value = 23
$ws.Range("AK4").Value = 21
$ws.Range("AN4").Value = 13.5

# Row 5
$ws.Range("G5").Value = 1.36
$ws.Range("N5").Value = 3.9
$ws.Range("P5").Value = 1.99
$ws.Range("Q5").Value = 1.96
$ws.Range("S5").Value = 3.5
$ws.Range("T5").Value = 2.46
$ws.Range("U5").Value = 1.64
$ws.Range("Z5").Value = 130
$ws.Range("AD5").Value = 50
$ws.Range("AE5").Value = 360
$ws.Range("AF5").Value = 7
$ws.Range("AN5").Value = 6.8

# Row 6
$ws.Range("H6").Value = 2.48
$ws.Range("I6").Value = 3.05
$ws.Range("Q6").Value = 1.82

# Row 7
$ws.Range("F7").Value = 1.57
$ws.Range("G7").Value = 1.71
$ws.Range("H7").Value = 5.4
$ws.Range("I7").Value = 7.4
$ws.Range("J7").Value = 4.2
$ws.Range("K7").Value = 5
$ws.Range("P7").Value = 2.16
$ws.Range("R7").Value = 1.43
$ws.Range("S7").Value = 2.6
$ws.Range("Y7").Value = 26
$ws.Range("AD7").Value = 27
$ws.Range("AL7").Value = 42
$ws.Range("AN7").Value = 10

# Row 8
$ws.Range("I8").Value = 1.29
$ws.Range("K8").Value = 7.8
$ws.Range("P8").Value = 4.2
$ws.Range("S8").Value = 1.73
$ws.Range("U8").Value = 2.5
$ws.Range("AA8").Value = 13.5
$ws.Range("AC8").Value = 20
$ws.Range("AH8").Value = 27
$ws.Range("AJ8").Value = 390
$ws.Range("AL8").Value = 85
$ws.Range("AM8").Value = 85
$ws.Range("AN8").Value = 100
$ws.Range("AO8").Value = 2.94

# Row 9
$ws.Range("G9").Value = 1.93
$ws.Range("H9").Value = 4.1
$ws.Range("I9").Value = 4.2
$ws.Range("Q9").Value = 1.65
$ws.Range("S9").Value = 2.6
$ws.Range("AA9").Value = 75
$ws.Range("AD9").Value = 17
$ws.Range("AE9").Value = 42
$ws.Range("AF9").Value = 14.5

# Row 10
$ws.Range("F10").Value = 1.8
$ws.Range("G10").Value = 1.82
$ws.Range("H10").Value = 4.7
$ws.Range("I10").Value = 4.9
$ws.Range("S10").Value = 2.62
$ws.Range("U10").Value = 2.38
$ws.Range("AB10").Value = 12
$ws.Range("AN10").Value = 8.4
$ws.Range("AO10").Value = 44

# Row 11
$ws.Range("K11").Value = 10
$ws.Range("P11").Value = 2.7
$ws.Range("Q11").Value = 1.57
$ws.Range("S11").Value = 2.38
$ws.Range("AD11").Value = 100
$ws.Range("AG11").Value = 14
$ws.Range("AH11").Value = 60
$ws.Range("AI11").Value = 490
$ws.Range("AJ11").Value = 7.4
$ws.Range("AM11").Value = 500

# Row 12
$ws.Range("F12").Value = 1.44
$ws.Range("G12").Value = 1.46
$ws.Range("H12").Value = 8.8
$ws.Range("I12").Value = 9.4

# Row 13
$ws.Range("F13").Value = 6.8
$ws.Range("G13").Value = 7.8
$ws.Range("H13").Value = 1.52
$ws.Range("J13").Value = 4.7
$ws.Range("O13").Value = 1.2
$ws.Range("Q13").Value = 1.6
$ws.Range("T13").Value = 1.79

# Row 14
$ws.Range("F14").Value = 1.23
$ws.Range("G14").Value = 1.5
$ws.Range("H14").Value = 1.04
$ws.Range("I14").Value = 16
$ws.Range("J14").Value = 4.5
$ws.Range("K14").Value = 980
$ws.Range("M14").Value = 1.01
$ws.Range("N14").Value = 1.04
$ws.Range("O14").Value = 1.2
$ws.Range("P14").Value = 1.5
$ws.Range("Q14").Value = 1.2
$ws.Range("R14").Value = 1.5
$ws.Range("S14").Value = 2.3
$ws.Range("T14").Value = 1.01
$ws.Range("U14").Value = 1.89
$ws.Range("V14").Value = 1.06
$ws.Range("W14").Value = 2.96
$ws.Range("X14").Value = 1000
$ws.Range("Y14").Value = 1000
$ws.Range("Z14").Value = 1000
$ws.Range("AA14").Value = 1000
$ws.Range("AB14").Value = 1000
$ws.Range("AC14").Value = 1000
$ws.Range("AD14").Value = 1000
$ws.Range("AE14").Value = 1000
$ws.Range("AF14").Value = 1000
$ws.Range("AG14").Value = 1000
$ws.Range("AH14").Value = 1000
$ws.Range("AI14").Value = 1000
$ws.Range("AJ14").Value = 1000
$ws.Range("AK14").Value = 1000
$ws.Range("AL14").Value = 1000
$ws.Range("AM14").Value = 1000
$ws.Range("AN14").Value = 1000
